# Generate Report for Handoff
# Replace the old handback-request GUID/hash with the new one that was
# generated for this handoff cycle, refresh the generation timestamps, and
# clear out the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns on the per-locale sheets (no handback has
# happened yet for the newly generated files).

$wb = $excel.ActiveWorkbook

$oldGuid = "b899ee7f-b814-4f33-b06f-7c25ceb5129d"
$newGuid = "1da57008-2498-4886-b7e0-5e4fb6b00348"
$oldHash = "bca88a58e47b48dfe857b994c05e7f83def8fddd"
$newHash = "3080efaf7b0a3c7eced77f984c433a9dac1699d1"

$oldGenDate = "2016-08-19 06:57:45"
$newGenDate = "2016-08-19 06:58:07"

$epoch = "0001-01-01 00:00:00"

# Cornflower-blue (FF6495ED), matching the workbook's original custom
# "HyperLink" cell style -- expressed as an OLE BGR color value for the
# Font.Color property.
$hyperlinkColor = 15570276

function Set-HyperlinkLook($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("B2").Value = "e2e\$newGuid.md"
$ws1.Range("G2").Value = $newGenDate

# Keep the same link target, just refresh the display text shown for it.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48cd87a385c15cf97d585a2460eab1204a1ca57d/e2e/$oldGuid.md", "", "", "e2e\$newGuid.md") | Out-Null
Set-HyperlinkLook $ws1.Range("B2")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-19 06:57:57"
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = $epoch

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48cd87a385c15cf97d585a2460eab1204a1ca57d/e2e/$oldGuid.md", "", "", "$newGuid.md") | Out-Null
Set-HyperlinkLook $ws2.Range("A2")

$ws2.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws2.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("H2").Value = $newGenDate
$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = $epoch

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48cd87a385c15cf97d585a2460eab1204a1ca57d/e2e/$oldGuid.md", "", "", "$newGuid.md") | Out-Null
Set-HyperlinkLook $ws3.Range("A2")

$ws3.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws3.Columns.Item(10).ColumnWidth = 21.7054770333426
